$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.127.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.85%  "
$ws.Range("D3").Value = "'2.249.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.22%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'245.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.92%  "
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").Value = "'75.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.42%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +7.60%  "
$ws.Range("D10").Value = "'40.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.88%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("E12").Value = "  +5.16%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'2.586.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "'14.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").Value = "'2.262.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "'0.798"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "'43.040.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("E19").Value = "  +6.06%  "
$ws.Range("D20").Value = "'71.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'5.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.12%  "
$ws.Range("D22").Value = "'9.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.09%  "
$ws.Range("D23").Value = "'230.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.20%  "
$ws.Range("E24").Value = "  +17.05%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'10.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("D27").Value = "'3.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").Value = "'39.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +30.29%  "
$ws.Range("D30").Value = "'2.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'173.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("D32").Value = "'20.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("D33").Value = "'0.0798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.35%  "
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").Value = "'0.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.47%  "
$ws.Range("E37").Value = "  +7.99%  "
$ws.Range("D38").Value = "'0.0335"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +20.55%  "
$ws.Range("D39").Value = "'13.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.23%  "
$ws.Range("E40").Value = "  +4.52%  "
$ws.Range("E41").Value = "  +3.66%  "
$ws.Range("D42").Value = "'0.206"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.75%  "
$ws.Range("D43").Value = "'59.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").Value = "'105.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.61%  "
$ws.Range("D45").Value = "'8.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.59%  "
$ws.Range("D46").Value = "'0.485"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +31.66%  "
$ws.Range("D47").Value = "'0.0995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.41%  "
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Value = "'2.460.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.85%  "
